# Auto-generated edit script for REV-11-03-2025.xlsx update
# - clears stray empty trailing cells on row 23
# - appends rows 24-53 with new product review data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23 cleanup: drop the empty trailing cells ---
$ws.Range("I23").ClearContents()
$ws.Range("K23:R23").ClearContents()

# --- Row 24 ---
$ws.Cells.Item(24, 1).Value = '0TF26835'
$ws.Cells.Item(24, 2).Value = 'I-AM UNAU BAKUCHIOL VANILLA REWID MASK'
$ws.Cells.Item(24, 3).Value = 'Consumo'
$ws.Cells.Item(24, 4).Value = 'No Tiene PT - TRADUZIDO'
$ws.Cells.Item(24, 5).Value = 'No Tiene ES - TRADUCIDO'
$ws.Cells.Item(24, 6).Value = 'No Tiene IT - TRADOTTO'
$ws.Cells.Item(24, 7).Value = '27'
$ws.Cells.Item(24, 8).Value = 'ML'
$ws.Cells.Item(24, 10).Value = 'Revisado y Traducido'
$ws.Cells.Item(24, 15).Value = 'Maschera Viso |Restauratrice| 27ml'
$ws.Cells.Item(24, 16).Value = 'Effetto: Aiuta a schiarire la pelle e migliora le rughe del 
viso.
Dopo aver deterso il viso, uniforma la texture della pelle 
con un tonico e applica la maschera sul viso. Lascia agire 
per 15-20 minuti, rimuovi la maschera e massaggia l''essenza 
residua.'
$ws.Cells.Item(24, 17).Value = 'Solo per uso esterno. Conservare in un luogo fresco e 
asciutto.'

# --- Row 25 ---
$ws.Cells.Item(25, 1).Value = '0TN03640'
$ws.Cells.Item(25, 2).Value = 'MONSTER HIGH CLAWDEEN WOLF HAND CREAM'
$ws.Cells.Item(25, 3).Value = 'Consumo'
$ws.Cells.Item(25, 4).Value = 'No Tiene PT - TRADUZIDO'
$ws.Cells.Item(25, 5).Value = 'Tiene ES'
$ws.Cells.Item(25, 6).Value = 'No Tiene IT - TRADOTTO'
$ws.Cells.Item(25, 7).Value = '30'
$ws.Cells.Item(25, 8).Value = 'GR'
$ws.Cells.Item(25, 10).Value = 'Revisado y Traducido'

# --- Row 26 ---
$ws.Cells.Item(26, 1).Value = 'ZMC00919'
$ws.Cells.Item(26, 2).Value = 'INSECT MATON 600 SP CASA Y JARDIN'
$ws.Cells.Item(26, 3).Value = 'Consumo'
$ws.Cells.Item(26, 4).Value = 'Tiene PT'
$ws.Cells.Item(26, 5).Value = 'Tiene ES'
$ws.Cells.Item(26, 6).Value = 'No Tiene IT - TRADOTTO'
$ws.Cells.Item(26, 7).Value = '600'
$ws.Cells.Item(26, 8).Value = 'ML'
$ws.Cells.Item(26, 10).Value = 'Revisado y Traducido'

# --- Row 27 ---
$ws.Cells.Item(27, 1).Value = '3AF01720'
$ws.Cells.Item(27, 2).Value = 'ENERGY FEELS ALWAYS LOVE EDT 28ML'
$ws.Cells.Item(27, 3).Value = 'Consumo'
$ws.Cells.Item(27, 4).Value = 'No Tiene PT - TRADUZIDO'
$ws.Cells.Item(27, 5).Value = 'Tiene ES'
$ws.Cells.Item(27, 6).Value = 'No Tiene IT - TRADOTTO'
$ws.Cells.Item(27, 7).Value = '28'
$ws.Cells.Item(27, 8).Value = 'ML'
$ws.Cells.Item(27, 10).Value = 'Revisado y Traducido'
$ws.Cells.Item(27, 11).Value = 'PTESTE'
$ws.Cells.Item(27, 12).Value = 'PTESTE'
$ws.Cells.Item(27, 13).Value = 'PTESTE'
$ws.Cells.Item(27, 14).Value = 'PTESTE'

# --- Row 28 ---
$ws.Cells.Item(28, 1).Value = '6VA40635'
$ws.Cells.Item(28, 2).Value = 'HELLO KITTY LIP BALM 7G BALL 3D'
$ws.Cells.Item(28, 3).Value = 'Consumo'
$ws.Cells.Item(28, 4).Value = 'Tiene PT'
$ws.Cells.Item(28, 5).Value = 'Tiene ES'
$ws.Cells.Item(28, 6).Value = 'Tiene IT'
$ws.Cells.Item(28, 7).Value = '7'
$ws.Cells.Item(28, 8).Value = 'GR'
$ws.Cells.Item(28, 10).Value = 'Solo Revisión'

# --- Row 29 ---
$ws.Cells.Item(29, 1).Value = '3AF01721'
$ws.Cells.Item(29, 2).Value = 'ENERGY FEELS PARTY HARD EDT 28ML'
$ws.Cells.Item(29, 3).Value = 'Consumo'
$ws.Cells.Item(29, 4).Value = 'Tiene PT'
$ws.Cells.Item(29, 5).Value = 'Tiene ES'
$ws.Cells.Item(29, 6).Value = 'Tiene IT'
$ws.Cells.Item(29, 7).Value = '28'
$ws.Cells.Item(29, 8).Value = 'ML'
$ws.Cells.Item(29, 10).Value = 'Solo Revisión'

# --- Row 30 ---
$ws.Cells.Item(30, 1).Value = '3AF01719'
$ws.Cells.Item(30, 2).Value = 'ENERGY FEELS LIVE BRIGHT EDT 28ML'
$ws.Cells.Item(30, 3).Value = 'Consumo'
$ws.Cells.Item(30, 4).Value = 'Tiene PT'
$ws.Cells.Item(30, 5).Value = 'Tiene ES'
$ws.Cells.Item(30, 6).Value = 'Tiene IT'
$ws.Cells.Item(30, 7).Value = '28'
$ws.Cells.Item(30, 8).Value = 'ML'
$ws.Cells.Item(30, 10).Value = 'Solo Revisión'

# --- Row 31 ---
$ws.Cells.Item(31, 1).Value = '2CA06584'
$ws.Cells.Item(31, 2).Value = 'ECO STYLE OLIVE OIL GEL 236 ML + 50% GRATIS'
$ws.Cells.Item(31, 3).Value = 'Consumo'
$ws.Cells.Item(31, 4).Value = 'No Tiene PT - TRADUZIDO'
$ws.Cells.Item(31, 5).Value = 'Tiene ES'
$ws.Cells.Item(31, 6).Value = 'No Tiene IT - TRADOTTO'
$ws.Cells.Item(31, 7).Value = '340'
$ws.Cells.Item(31, 8).Value = 'ML'
$ws.Cells.Item(31, 10).Value = 'Revisado y Traducido'

# --- Row 32 ---
$ws.Cells.Item(32, 1).Value = '0TP03790'
$ws.Cells.Item(32, 2).Value = 'MY PRIDE SET BROCHAS'
$ws.Cells.Item(32, 3).Value = 'Consumo'
$ws.Cells.Item(32, 4).Value = 'Tiene PT'
$ws.Cells.Item(32, 5).Value = 'Tiene ES'
$ws.Cells.Item(32, 6).Value = 'Tiene IT'
$ws.Cells.Item(32, 7).Value = '1'
$ws.Cells.Item(32, 8).Value = 'UND'
$ws.Cells.Item(32, 10).Value = 'Solo Revisión'

# --- Row 33 ---
$ws.Cells.Item(33, 1).Value = '6VA38727'
$ws.Cells.Item(33, 2).Value = 'BILLY BROWN TIES EXECUTIVE SOCKS  36-40 M'
$ws.Cells.Item(33, 3).Value = 'Consumo'
$ws.Cells.Item(33, 4).Value = 'Tiene PT'
$ws.Cells.Item(33, 5).Value = 'Tiene ES'
$ws.Cells.Item(33, 6).Value = 'Tiene IT'
$ws.Cells.Item(33, 7).Value = '1'
$ws.Cells.Item(33, 8).Value = 'UND'
$ws.Cells.Item(33, 10).Value = 'Solo Revisión'

# --- Row 34 ---
$ws.Cells.Item(34, 1).Value = '6VA38734'
$ws.Cells.Item(34, 2).Value = 'BILLY BROWN APEROL ITALIAN SOCKS  41-46 H'
$ws.Cells.Item(34, 3).Value = 'Consumo'
$ws.Cells.Item(34, 4).Value = 'Tiene PT'
$ws.Cells.Item(34, 5).Value = 'Tiene ES'
$ws.Cells.Item(34, 6).Value = 'Tiene IT'
$ws.Cells.Item(34, 7).Value = '1'
$ws.Cells.Item(34, 8).Value = 'UND'
$ws.Cells.Item(34, 10).Value = 'Solo Revisión'

# --- Row 35 ---
$ws.Cells.Item(35, 1).Value = '6VA35313'
$ws.Cells.Item(35, 2).Value = 'INVISIBOBBLE HAIRHALO MARGARITA BONITA'
$ws.Cells.Item(35, 3).Value = 'Consumo'
$ws.Cells.Item(35, 4).Value = 'Tiene PT'
$ws.Cells.Item(35, 5).Value = 'Tiene ES'
$ws.Cells.Item(35, 6).Value = 'Tiene IT'
$ws.Cells.Item(35, 7).Value = '1'
$ws.Cells.Item(35, 8).Value = 'UND'
$ws.Cells.Item(35, 10).Value = 'Solo Revisión'

# --- Row 36 ---
$ws.Cells.Item(36, 1).Value = 'MASTER TESTE'
$ws.Cells.Item(36, 2).Value = 'MASTER TESTE'
$ws.Cells.Item(36, 3).Value = 'Consumo'
$ws.Cells.Item(36, 4).Value = 'No Tiene PT - TRADUZIDO'
$ws.Cells.Item(36, 5).Value = 'Tiene ES'
$ws.Cells.Item(36, 6).Value = 'No Tiene IT - TRADOTTO'
$ws.Cells.Item(36, 7).Value = '1'
$ws.Cells.Item(36, 8).Value = 'UND'
$ws.Cells.Item(36, 10).Value = 'Revisado y Traducido'
$ws.Cells.Item(36, 11).Value = 'MASTER TESTE PT TITULO'
$ws.Cells.Item(36, 12).Value = 'MASTER TESTE PT USO'
$ws.Cells.Item(36, 13).Value = 'MASTER TESTE PT ADVERTENCIAS'
$ws.Cells.Item(36, 14).Value = 'MASTER TESTE PT +INFO'
$ws.Cells.Item(36, 15).Value = 'MASTER TESTE IT TITULO'
$ws.Cells.Item(36, 16).Value = 'MASTER TESTE IT USO'
$ws.Cells.Item(36, 17).Value = 'MASTER TESTE IT ADVERTENCIAS'
$ws.Cells.Item(36, 18).Value = 'MASTER TESTE IT +INFO'

# --- Row 37 ---
$ws.Cells.Item(37, 1).Value = '0TP03782'
$ws.Cells.Item(37, 2).Value = 'W-7 LIP CARE KIT - VANILLA'
$ws.Cells.Item(37, 3).Value = 'Consumo'
$ws.Cells.Item(37, 4).Value = 'No Tiene PT - TRADUZIDO'
$ws.Cells.Item(37, 5).Value = 'Tiene ES'
$ws.Cells.Item(37, 6).Value = 'No Tiene IT - TRADOTTO'
$ws.Cells.Item(37, 7).Value = '2'
$ws.Cells.Item(37, 8).Value = 'UND'
$ws.Cells.Item(37, 10).Value = 'Revisado y Traducido'

# --- Row 38 ---
$ws.Cells.Item(38, 1).Value = '0TP03780'
$ws.Cells.Item(38, 2).Value = 'W-7 LIP CARE KIT - CHERRY'
$ws.Cells.Item(38, 3).Value = 'Consumo'
$ws.Cells.Item(38, 4).Value = 'No Tiene PT - TRADUZIDO'
$ws.Cells.Item(38, 5).Value = 'Tiene ES'
$ws.Cells.Item(38, 6).Value = 'No Tiene IT - TRADOTTO'
$ws.Cells.Item(38, 7).Value = '2'
$ws.Cells.Item(38, 8).Value = 'UND'
$ws.Cells.Item(38, 10).Value = 'Revisado y Traducido'

# --- Row 39 ---
$ws.Cells.Item(39, 1).Value = '1DT00927'
$ws.Cells.Item(39, 2).Value = 'THE FRUIT COMPANY PASTA DENTAL MELON 60GR'
$ws.Cells.Item(39, 3).Value = 'Consumo'
$ws.Cells.Item(39, 4).Value = 'Tiene PT'
$ws.Cells.Item(39, 5).Value = 'Tiene ES'
$ws.Cells.Item(39, 6).Value = 'Tiene IT'
$ws.Cells.Item(39, 7).Value = '60'
$ws.Cells.Item(39, 8).Value = 'GR'
$ws.Cells.Item(39, 10).Value = 'Solo Revisión'

# --- Row 40 ---
$ws.Cells.Item(40, 1).Value = '1DT00926'
$ws.Cells.Item(40, 2).Value = 'THE FRUIT COMPANY PASTA DENTAL SANDIA 60GR'
$ws.Cells.Item(40, 3).Value = 'Consumo'
$ws.Cells.Item(40, 4).Value = 'Tiene PT'
$ws.Cells.Item(40, 5).Value = 'Tiene ES'
$ws.Cells.Item(40, 6).Value = 'Tiene IT'
$ws.Cells.Item(40, 7).Value = '60'
$ws.Cells.Item(40, 8).Value = 'GR'
$ws.Cells.Item(40, 10).Value = 'Solo Revisión'

# --- Row 41 ---
$ws.Cells.Item(41, 1).Value = '1DT00929'
$ws.Cells.Item(41, 2).Value = 'THE FRUIT COMPANY PASTA DENTAL MELOCOTON 60GR'
$ws.Cells.Item(41, 3).Value = 'Consumo'
$ws.Cells.Item(41, 4).Value = 'Tiene PT'
$ws.Cells.Item(41, 5).Value = 'Tiene ES'
$ws.Cells.Item(41, 6).Value = 'Tiene IT'
$ws.Cells.Item(41, 7).Value = '60'
$ws.Cells.Item(41, 8).Value = 'GR'
$ws.Cells.Item(41, 10).Value = 'Solo Revisión'

# --- Row 42 ---
$ws.Cells.Item(42, 1).Value = '3AF01926'
$ws.Cells.Item(42, 2).Value = 'SJP LOVELY LOVELY LIGHTS BODY MIST 236ML'
$ws.Cells.Item(42, 3).Value = 'Consumo'
$ws.Cells.Item(42, 4).Value = 'Tiene PT'
$ws.Cells.Item(42, 5).Value = 'Tiene ES'
$ws.Cells.Item(42, 6).Value = 'Tiene IT'
$ws.Cells.Item(42, 7).Value = '236'
$ws.Cells.Item(42, 8).Value = 'ML'
$ws.Cells.Item(42, 10).Value = 'Solo Revisión'

# --- Row 43 ---
$ws.Cells.Item(43, 1).Value = '3AF01938'
$ws.Cells.Item(43, 2).Value = 'TAI & JON BODYMIST DOLCE ICE CREAM 250ML'
$ws.Cells.Item(43, 3).Value = 'Consumo'
$ws.Cells.Item(43, 4).Value = 'Tiene PT'
$ws.Cells.Item(43, 5).Value = 'Tiene ES'
$ws.Cells.Item(43, 6).Value = 'Tiene IT'
$ws.Cells.Item(43, 7).Value = '250'
$ws.Cells.Item(43, 8).Value = 'ML'
$ws.Cells.Item(43, 10).Value = 'Solo Revisión'

# --- Row 44 ---
$ws.Cells.Item(44, 1).Value = '3AF01940'
$ws.Cells.Item(44, 2).Value = 'TAI & JON BODYMIST FRUIT TRANQUILITY 250ML'
$ws.Cells.Item(44, 3).Value = 'Consumo'
$ws.Cells.Item(44, 4).Value = 'Tiene PT'
$ws.Cells.Item(44, 5).Value = 'Tiene ES'
$ws.Cells.Item(44, 6).Value = 'Tiene IT'
$ws.Cells.Item(44, 7).Value = '250'
$ws.Cells.Item(44, 8).Value = 'ML'
$ws.Cells.Item(44, 10).Value = 'Solo Revisión'

# --- Row 45 ---
$ws.Cells.Item(45, 1).Value = '6XS18308'
$ws.Cells.Item(45, 2).Value = 'CERAVE CREMA HIDRATANTE REFILL 454GRS'
$ws.Cells.Item(45, 3).Value = 'Consumo'
$ws.Cells.Item(45, 4).Value = 'Tiene PT'
$ws.Cells.Item(45, 5).Value = 'Tiene ES'
$ws.Cells.Item(45, 6).Value = 'Tiene IT'
$ws.Cells.Item(45, 7).Value = '454'
$ws.Cells.Item(45, 8).Value = 'GR'
$ws.Cells.Item(45, 10).Value = 'Solo Revisión'

# --- Row 46 ---
$ws.Cells.Item(46, 1).Value = '6XS18295'
$ws.Cells.Item(46, 2).Value = 'CERAVE LOCION HIDRATANTE REFILL 473ML'
$ws.Cells.Item(46, 3).Value = 'Consumo'
$ws.Cells.Item(46, 4).Value = 'Tiene PT'
$ws.Cells.Item(46, 5).Value = 'Tiene ES'
$ws.Cells.Item(46, 6).Value = 'Tiene IT'
$ws.Cells.Item(46, 7).Value = '473'
$ws.Cells.Item(46, 8).Value = 'ML'
$ws.Cells.Item(46, 10).Value = 'Solo Revisión'

# --- Row 47 ---
$ws.Cells.Item(47, 1).Value = '6XS18293'
$ws.Cells.Item(47, 2).Value = 'CERAVE GEL LIMPIADOR ESPUMOSO REFILL 473ML'
$ws.Cells.Item(47, 3).Value = 'Consumo'
$ws.Cells.Item(47, 4).Value = 'Tiene PT'
$ws.Cells.Item(47, 5).Value = 'Tiene ES'
$ws.Cells.Item(47, 6).Value = 'Tiene IT'
$ws.Cells.Item(47, 7).Value = '473'
$ws.Cells.Item(47, 8).Value = 'ML'
$ws.Cells.Item(47, 10).Value = 'Solo Revisión'

# --- Row 48 ---
$ws.Cells.Item(48, 1).Value = '6XS18294'
$ws.Cells.Item(48, 2).Value = 'CERAVE LIMPIADORA HIDRATANTE REFILL 473ML'
$ws.Cells.Item(48, 3).Value = 'Consumo'
$ws.Cells.Item(48, 4).Value = 'Tiene PT'
$ws.Cells.Item(48, 5).Value = 'Tiene ES'
$ws.Cells.Item(48, 6).Value = 'Tiene IT'
$ws.Cells.Item(48, 7).Value = '473'
$ws.Cells.Item(48, 8).Value = 'ML'
$ws.Cells.Item(48, 10).Value = 'Solo Revisión'

# --- Row 49 ---
$ws.Cells.Item(49, 1).Value = '2LT03571'
$ws.Cells.Item(49, 2).Value = 'MET EXPOSITOR GOTAS BRONCEADORAS'
$ws.Cells.Item(49, 3).Value = 'LOTE'
$ws.Cells.Item(49, 4).Value = 'Tiene PT'
$ws.Cells.Item(49, 5).Value = 'Tiene ES'
$ws.Cells.Item(49, 6).Value = 'Tiene IT'
$ws.Cells.Item(49, 7).Value = '10'
$ws.Cells.Item(49, 8).Value = 'UND'
$ws.Cells.Item(49, 9).Value = '"8445984023355"'
$ws.Cells.Item(49, 10).Value = 'Solo Revisión'

# --- Row 50 ---
$ws.Cells.Item(50, 1).Value = '0MR27448'
$ws.Cells.Item(50, 2).Value = 'MET BRONZE DROPS GOTAS BRONCEADORAS 30ML'
$ws.Cells.Item(50, 3).Value = 'Consumo'
$ws.Cells.Item(50, 4).Value = 'No Tiene PT - TRADUZIDO'
$ws.Cells.Item(50, 5).Value = 'Tiene ES'
$ws.Cells.Item(50, 6).Value = 'No Tiene IT - TRADOTTO'
$ws.Cells.Item(50, 7).Value = '30'
$ws.Cells.Item(50, 8).Value = 'ML'
$ws.Cells.Item(50, 10).Value = 'Revisado y Traducido'

# --- Row 51 ---
$ws.Cells.Item(51, 1).Value = '2LT03460'
$ws.Cells.Item(51, 2).Value = 'TECHNIC LOTE JELLY BLUSH - PINK BURST'
$ws.Cells.Item(51, 3).Value = 'LOTE'
$ws.Cells.Item(51, 4).Value = 'Tiene PT'
$ws.Cells.Item(51, 5).Value = 'Tiene ES'
$ws.Cells.Item(51, 6).Value = 'Tiene IT'
$ws.Cells.Item(51, 7).Value = '10'
$ws.Cells.Item(51, 8).Value = 'UND'
$ws.Cells.Item(51, 9).Value = '"5021769247371"'
$ws.Cells.Item(51, 10).Value = 'Solo Revisión'

# --- Row 52 ---
$ws.Cells.Item(52, 1).Value = '0MR27166'
$ws.Cells.Item(52, 2).Value = 'TECHNIC JELLY BLUSH - PINK BURST'
$ws.Cells.Item(52, 3).Value = 'Consumo'
$ws.Cells.Item(52, 4).Value = 'Tiene PT'
$ws.Cells.Item(52, 5).Value = 'Tiene ES'
$ws.Cells.Item(52, 6).Value = 'Tiene IT'
$ws.Cells.Item(52, 7).Value = '9'
$ws.Cells.Item(52, 8).Value = 'GR'
$ws.Cells.Item(52, 10).Value = 'Solo Revisión'

# --- Row 53 ---
$ws.Cells.Item(53, 1).Value = 'teste sku'
$ws.Cells.Item(53, 2).Value = 'teste titulo'
$ws.Cells.Item(53, 3).Value = 'Consumo'
$ws.Cells.Item(53, 4).Value = 'No Tiene PT - TRADUZIDO'
$ws.Cells.Item(53, 5).Value = 'Tiene ES'
$ws.Cells.Item(53, 6).Value = 'No Tiene IT - TRADOTTO'
$ws.Cells.Item(53, 7).Value = '1'
$ws.Cells.Item(53, 8).Value = 'UND'
$ws.Cells.Item(53, 10).Value = 'Revisado y Traducido'
$ws.Cells.Item(53, 11).Value = 'teste'
$ws.Cells.Item(53, 12).Value = 'teste'
$ws.Cells.Item(53, 13).Value = 'teste'
$ws.Cells.Item(53, 14).Value = 'teste'
$ws.Cells.Item(53, 15).Value = 'teste1'
$ws.Cells.Item(53, 16).Value = 'teste1'
$ws.Cells.Item(53, 17).Value = 'teste1'
$ws.Cells.Item(53, 18).Value = 'teste1'

